$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny correction to the existing timestamp in A8 (automatic re-sync from WSL)
$ws.Range("A8").Value = 45878.29184623843

# Append the new hourly reading as row 9
$ws.Range("A9").Value = 45878.33351998493
$ws.Range("B9").Value = 2025
$ws.Range("C9").Value = 37
$ws.Range("D9").Value = 13.28
$ws.Range("E9").Value = 92.75
$ws.Range("F9").Value = 77.34
$ws.Range("G9").Value = 7.87
$ws.Range("H9").Value = "ESE"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "08:00:16"

# Match the date/time style used by A2:A8
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
